$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.530.64"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "1.736.40"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("E7").Value = "  +1.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06225"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("D10").Value = "1.730.41"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07035"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.596"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6098"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.49"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "26.518.89"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("E18").Value = "  +6.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").Value = "  -2.15%  "

$ws.Range("D21").Value = "1.954.70"
$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.555"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.766"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.236"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.72"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.56%  "

$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.414"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.775"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "108.08"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.012"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08069"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.694"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04563"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9995"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.611"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.009"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6359"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9012"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.027"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.397"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01505"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.51"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.425"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.74%  "

$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.937"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05393"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.12%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.818"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.87%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.257"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.03%  "
